# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Femacal de La Calera - Choclo" right
# before the existing row 544, pushing the old rows 544-553 down to 547-556.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 544 (old data shifts down).
$ws.Range("A544:A546").EntireRow.Insert()

# New weekly data (columns A..R) for the freshly-inserted rows 544-546.
$newData = @(
    @(3, "Femacal de La Calera", "Coquimbo", 44595, 5, 100112024, "Choclo", "Choclero", "Primera", 17000, 200, 250, 225, "`$/unidad", "Provincia de Quillota", 225, 1, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44595, 5, 100112024, "Choclo", "Choclero", "Segunda", 7000, 150, 150, 150, "`$/unidad", "Provincia de Quillota", 150, 1, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44595, 5, 100112024, "Choclo", "Dulce o Americano", "Primera", 15100, 150, 180, 165, "`$/unidad", "Provincia de Quillota", 165, 1, "Hortaliza")
)

$startRow = 544
for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowNum = $startRow + $i
    $rowVals = $newData[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
